# mpsi2_TP.xlsx : maj info 2 premier chapitre et TP
# Adds two new columns (M = "DS9", N = "Projet") to the "S2" sheet, mirroring
# the shading pattern of the existing table, plus a totals formula in M46,
# and updates the window/selection view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (M1/N1), styled like the existing L1 "DS7" header cell ----
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 13).Value = "DS9"
$ws.Cells.Item(1, 14).Value = "Projet"

# ---- Per-row data, columns M (13) and N (14) ----
# "style"  -> shaded cell matching the row's fill (borrowed from A2, fillId=2,
#             border removed) == cellXf used by rows 4,8,10,12,... in the diff
# "plain"  -> value written with the default/general style (no explicit xf)
# $null    -> cell intentionally left absent for that row
$rowData = @(
    @{Row=2; M="style"; N="plain"},
    @{Row=3; M="plain"; N="plain"},
    @{Row=4; M="style"; N="style"},
    @{Row=5; M="plain"; N=$null},
    @{Row=6; M="plain"; N="plain"},
    @{Row=7; M="plain"; N="plain"},
    @{Row=8; M="style"; N="style"},
    @{Row=9; M="plain"; N="plain"},
    @{Row=10; M="style"; N="style"},
    @{Row=11; M="plain"; N="plain"},
    @{Row=12; M="style"; N="style"},
    @{Row=13; M="plain"; N="plain"},
    @{Row=14; M="style"; N="style"},
    @{Row=15; M="plain"; N="plain"},
    @{Row=16; M="style"; N="style"},
    @{Row=17; M="plain"; N="plain"},
    @{Row=18; M="style"; N="style"},
    @{Row=19; M="plain"; N="plain"},
    @{Row=20; M="style"; N="style"},
    @{Row=21; M="plain"; N="plain"},
    @{Row=22; M="style"; N="style"},
    @{Row=23; M="plain"; N="plain"},
    @{Row=24; M="style"; N="style"},
    @{Row=25; M="plain"; N="plain"},
    @{Row=26; M="style"; N="plain"},
    @{Row=27; M="plain"; N=$null},
    @{Row=28; M="style"; N="plain"},
    @{Row=29; M="plain"; N="plain"},
    @{Row=30; M="style"; N="style"},
    @{Row=31; M="plain"; N="plain"},
    @{Row=32; M="plain"; N="plain"},
    @{Row=33; M="style"; N="style"},
    @{Row=34; M="plain"; N="plain"},
    @{Row=35; M="style"; N="style"},
    @{Row=36; M="plain"; N="plain"},
    @{Row=37; M="plain"; N="plain"},
    @{Row=38; M="style"; N="style"},
    @{Row=39; M="style"; N="style"},
    @{Row=40; M=$null; N="plain"},
    @{Row=41; M="style"; N="style"},
    @{Row=42; M="plain"; N="plain"},
    @{Row=43; M="style"; N="style"},
    @{Row=44; M="plain"; N="plain"},
    @{Row=45; M=$null; N=$null}
)

foreach ($item in $rowData) {
    $r = $item.Row

    if ($item.M -eq "style") {
        $ws.Range("A2").Copy()
        $ws.Cells.Item($r, 13).PasteSpecial(-4122)
        $ws.Application.CutCopyMode = $false
        $ws.Cells.Item($r, 13).Borders.LineStyle = -4142
        $ws.Cells.Item($r, 13).Value = 1
    } elseif ($item.M -eq "plain") {
        $ws.Cells.Item($r, 13).Value = 1
    }

    if ($item.N -eq "style") {
        $ws.Range("A2").Copy()
        $ws.Cells.Item($r, 14).PasteSpecial(-4122)
        $ws.Application.CutCopyMode = $false
        $ws.Cells.Item($r, 14).Borders.LineStyle = -4142
        $ws.Cells.Item($r, 14).Value = 1
    } elseif ($item.N -eq "plain") {
        $ws.Cells.Item($r, 14).Value = 1
    }
}

# ---- Totals row ----
$ws.Cells.Item(46, 13).Formula = "=SUM(M2:M45)"

# ---- View state: selection + scroll position ----
$ws.Activate() | Out-Null
$excel.Goto($ws.Range("A19"), $true) | Out-Null
$ws.Range("N50").Select() | Out-Null

# ---- Workbook window geometry ----
try { $excel.Left = 2380 } catch {}
try { $excel.Top = 1560 } catch {}
try { $excel.Width = 19400 } catch {}
try { $excel.Height = 14220 } catch {}
